{"js": "// \"Trocando tags no template\" \u2014 remove the stray \"-7\" suffix from the\n// bank-account number (\"Conta Banc\u00e1ria: 10.738-7\" -> \"... 10.738\") and\n// remove the hard-coded validity date range after \"Vig\u00eancia das\n// Atividades do Projeto:\" (the field becomes blank, ready to be filled\n// in by the template engine).\n\nconst body = context.document.body;\n\n// 1) \"Conta Banc\u00e1ria: 10.738-7\" -> \"Conta Banc\u00e1ria: 10.738\"\nconst accountSuffix = body.search(\"-7\", { matchCase: true, matchWholeWord: false });\naccountSuffix.load(\"items/text\");\nawait context.sync();\n\nfor (const r of accountSuffix.items) {\n  if (r.text === \"-7\") {\n    r.delete();\n  }\n}\nawait context.sync();\n\n// 2) \"Vig\u00eancia das Atividades do Projeto: 23/02/2021 a 23/02/2022\" ->\n//    \"Vig\u00eancia das Atividades do Projeto:\"\nconst vigenciaDates = body.search(\" 23/02/2021 a 23/02/2022\", { matchCase: true, matchWholeWord: false });\nvigenciaDates.load(\"items/text\");\nawait context.sync();\n\nfor (const r of vigenciaDates.items) {\n  r.delete();\n}\nawait context.sync();\n", "ps1": "# \"Trocando tags no template\" \u2014 remove the stray \"-7\" suffix from the\n# bank-account number (\"Conta Bancaria: 10.738-7\" -> \"... 10.738\") and\n# remove the hard-coded validity date range after \"Vigencia das\n# Atividades do Projeto:\" (the field becomes blank, ready to be filled\n# in by the template engine).\n\n$d = $word.ActiveDocument\n\n# 1) \"Conta Banc\u00e1ria: 10.738-7\" -> \"Conta Banc\u00e1ria: 10.738\"\n$r1 = $d.Content\n$found1 = $r1.Find.Execute(\"-7\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found1) {\n    $r1.Delete()\n}\n\n# 2) \"Vig\u00eancia das Atividades do Projeto: 23/02/2021 a 23/02/2022\" ->\n#    \"Vig\u00eancia das Atividades do Projeto:\"\n$r2 = $d.Content\n$found2 = $r2.Find.Execute(\" 23/02/2021 a 23/02/2022\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\nif ($found2) {\n    $r2.Delete()\n}\n"}
